$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "core" (student) field table: add required-field "*" markers in column I ---
$ws.Range("I3").Value = "*"
$ws.Range("I4").Value = "*"
$ws.Range("I5").Value = "*"
$ws.Range("I6").Value = "*"
$ws.Range("K9").Value = "*"
$ws.Range("I10").Value = "*"
$ws.Range("I11").Value = "*"
$ws.Range("J23").Value = "*"
$ws.Range("I26").Value = "*"

# --- "teacher" table: relocate from K:L to N:O, add "*" markers in column P ---
$ws.Range("K3").Value = $null
$ws.Range("L3").Value = $null
$ws.Range("L4").Value = $null
$ws.Range("L5").Value = $null
$ws.Range("L6").Value = $null
$ws.Range("L7").Value = $null
$ws.Range("L8").Value = $null
$ws.Range("L9").Value = $null
$ws.Range("L10").Value = $null
$ws.Range("L11").Value = $null

$ws.Range("N3").Value = "teacher"
$ws.Range("O3").Value = "id"
$ws.Range("P3").Value = "*"
$ws.Range("O4").Value = "firstName"
$ws.Range("P4").Value = "*"
$ws.Range("O5").Value = "lastName"
$ws.Range("P5").Value = "*"
$ws.Range("O6").Value = "code"
$ws.Range("P6").Value = "*"
$ws.Range("O7").Value = "department"
$ws.Range("P7").Value = "*"
$ws.Range("O8").Value = "gender"
$ws.Range("P8").Value = "*"
$ws.Range("O9").Value = "phone"
$ws.Range("P9").Value = "*"
$ws.Range("O10").Value = "dob"
$ws.Range("P10").Value = "*"
$ws.Range("O11").Value = "email"
$ws.Range("P11").Value = "*"

# --- "subject" table: relocate from N:O to Q:R, add "*" markers in column S ---
$ws.Range("Q3").Value = "subject"
$ws.Range("R3").Value = "id"
$ws.Range("S3").Value = "*"
$ws.Range("R4").Value = "name"
$ws.Range("S4").Value = "*"
$ws.Range("R5").Value = "code"
$ws.Range("S5").Value = "*"
$ws.Range("R6").Value = "department"
$ws.Range("S6").Value = "*"
$ws.Range("R7").Value = "all"
$ws.Range("S7").Value = "*"
$ws.Range("R8").Value = "theory"
$ws.Range("R9").Value = "practice"
$ws.Range("R10").Value = "exercise"

# --- selection, as left by the editing session ---
$ws.Range("B6").Select()
